# Periodic data refresh ("updated 4.0 files and mdl"):
#  - bump the "last updated" date on the About sheet
#  - update the hard coal production/import/export balancing priorities
#  - leave the user's selection on the FPIEBP sheet at E3

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsFP    = $wb.Worksheets.Item("FPIEBP")

# About!C1 - "last updated" date serial (2024-01-03 -> 2024-03-28)
$wsAbout.Range("C1").Value = 45379

# FPIEBP!B3:D3 - "hard coal" priorities: production/imports/exports
$wsFP.Range("B3").Value = 1
$wsFP.Range("C3").Value = 3
$wsFP.Range("D3").Value = 2

# Leave the FPIEBP sheet active with E3 selected
$wsFP.Activate()
$wsFP.Range("E3").Select()
